# "Changes in add HSIM menu"
# Update the sample login row: username changes from "BBSSL4" to "BBSSL92".
# (Password/result columns, and the hyperlink on B2, are left untouched.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "BBSSL92"

# Move/restore the active selection to A2 (matches the saved workbook view).
[void]$ws.Range("A2").Select()
